$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E3").Value = 32
$ws.Range("E5").Value = 33
$ws.Range("E6").Value = 60
$ws.Range("E12").Value = 35
$ws.Range("E15").Value = 114
$ws.Range("F16").Value = 97
$ws.Range("H16").Value = 185
$ws.Range("E18").Value = 100
